$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly refresh re-sorts/permutes the existing daily price records
# (date, quality, volume, min/max/avg price, price per kg) across rows
# 2-33, while the descriptive columns (market, product, category, unit,
# origin, etc.) stay identical for every row. Capture a snapshot of the
# columns that move (D, L, M, N, O, P, S) before writing the new layout
# so that the permutation can be applied safely regardless of cycle
# order.

$firstRow = 2
$lastRow = 33

$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = @(
        $ws.Cells.Item($r, 4).Value2,   # D - Fecha
        $ws.Cells.Item($r, 12).Value2,  # L - Calidad
        $ws.Cells.Item($r, 13).Value2,  # M - Volumen
        $ws.Cells.Item($r, 14).Value2,  # N - Precio minimo
        $ws.Cells.Item($r, 15).Value2,  # O - Precio maximo
        $ws.Cells.Item($r, 16).Value2,  # P - Precio promedio ponderado
        $ws.Cells.Item($r, 19).Value2   # S - Precio $/Kg
    )
}

# Maps each destination row to the row whose record it now receives.
$rowMap = @{
    2  = 26
    3  = 27
    4  = 28
    5  = 14
    6  = 15
    7  = 21
    8  = 22
    9  = 29
    10 = 32
    11 = 33
    12 = 2
    13 = 3
    14 = 4
    15 = 9
    16 = 10
    17 = 13
    18 = 11
    19 = 12
    20 = 18
    21 = 19
    22 = 20
    23 = 5
    24 = 6
    25 = 23
    26 = 24
    27 = 25
    28 = 16
    29 = 17
    30 = 7
    31 = 8
    32 = 30
    33 = 31
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]

    $ws.Cells.Item($destRow, 4).Value = $vals[0]
    $ws.Cells.Item($destRow, 12).Value = $vals[1]
    $ws.Cells.Item($destRow, 13).Value = $vals[2]
    $ws.Cells.Item($destRow, 14).Value = $vals[3]
    $ws.Cells.Item($destRow, 15).Value = $vals[4]
    $ws.Cells.Item($destRow, 16).Value = $vals[5]
    $ws.Cells.Item($destRow, 19).Value = $vals[6]
}
